# Atualizado por script em 25-10-2023 12:32
#
# 1) Rows 84 and 85 hold two matches that were played on the same day
#    (Mladost vs Radnicki Nis, and Vojvodina vs Radnik). The match-specific
#    columns (F..V, skipping the shared timestamp columns K/O/S which are
#    identical for both rows) are swapped between the two rows.
# 2) A brand new match (row 87 / Indice 86) is appended: Radnicki 1923 vs
#    Partizan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: swap the per-match data between row 84 and row 85.
# Columns A,B,C,D,E,K,O,S stay put; F,G,H,I,J,L,M,N,P,Q,R,T,U,V swap.
# ---------------------------------------------------------------------

$row84_F = $ws.Range("F84").Value2
$row84_G = $ws.Range("G84").Value2
$row84_H = $ws.Range("H84").Value2
$row84_I = $ws.Range("I84").Value2
$row84_J = $ws.Range("J84").Value2
$row84_L = $ws.Range("L84").Value2
$row84_M = $ws.Range("M84").Value2
$row84_N = $ws.Range("N84").Value2
$row84_P = $ws.Range("P84").Value2
$row84_Q = $ws.Range("Q84").Value2
$row84_R = $ws.Range("R84").Value2
$row84_T = $ws.Range("T84").Value2
$row84_U = $ws.Range("U84").Value2
$row84_V = $ws.Range("V84").Value2

$row85_F = $ws.Range("F85").Value2
$row85_G = $ws.Range("G85").Value2
$row85_H = $ws.Range("H85").Value2
$row85_I = $ws.Range("I85").Value2
$row85_J = $ws.Range("J85").Value2
$row85_L = $ws.Range("L85").Value2
$row85_M = $ws.Range("M85").Value2
$row85_N = $ws.Range("N85").Value2
$row85_P = $ws.Range("P85").Value2
$row85_Q = $ws.Range("Q85").Value2
$row85_R = $ws.Range("R85").Value2
$row85_T = $ws.Range("T85").Value2
$row85_U = $ws.Range("U85").Value2
$row85_V = $ws.Range("V85").Value2

$ws.Range("F84").Value2 = $row85_F
$ws.Range("G84").Value2 = $row85_G
$ws.Range("H84").Value2 = $row85_H
$ws.Range("I84").Value2 = $row85_I
$ws.Range("J84").Value2 = $row85_J
$ws.Range("L84").Value2 = $row85_L
$ws.Range("M84").Value2 = $row85_M
$ws.Range("N84").Value2 = $row85_N
$ws.Range("P84").Value2 = $row85_P
$ws.Range("Q84").Value2 = $row85_Q
$ws.Range("R84").Value2 = $row85_R
$ws.Range("T84").Value2 = $row85_T
$ws.Range("U84").Value2 = $row85_U
$ws.Range("V84").Value2 = $row85_V

$ws.Range("F85").Value2 = $row84_F
$ws.Range("G85").Value2 = $row84_G
$ws.Range("H85").Value2 = $row84_H
$ws.Range("I85").Value2 = $row84_I
$ws.Range("J85").Value2 = $row84_J
$ws.Range("L85").Value2 = $row84_L
$ws.Range("M85").Value2 = $row84_M
$ws.Range("N85").Value2 = $row84_N
$ws.Range("P85").Value2 = $row84_P
$ws.Range("Q85").Value2 = $row84_Q
$ws.Range("R85").Value2 = $row84_R
$ws.Range("T85").Value2 = $row84_T
$ws.Range("U85").Value2 = $row84_U
$ws.Range("V85").Value2 = $row84_V

# ---------------------------------------------------------------------
# Step 2: append the new match as row 87 (Indice 86).
# Copy formatting from row 86 (the current last row) first, then fill
# in the values.
# ---------------------------------------------------------------------

$ws.Range("A86:V86").Copy($ws.Range("A87:V87"))

$ws.Range("A87").Value2 = 86
$ws.Range("B87").Value2 = "serbia"
$ws.Range("C87").Value2 = "super-liga"
$ws.Range("D87").Value2 = "2023-2024"
$ws.Range("E87").Value2 = 45224.58333333334
$ws.Range("F87").Value2 = "Radnicki 1923"
$ws.Range("G87").Value2 = 0
$ws.Range("H87").Value2 = "Partizan"
$ws.Range("I87").Value2 = 4
$ws.Range("J87").Value2 = 5.64
$ws.Range("K87").Value2 = "17/08/2023 08:13"
$ws.Range("L87").Value2 = 9.619999999999999
$ws.Range("M87").Value2 = "25/10/2023 13:46"
$ws.Range("N87").Value2 = 4.04
$ws.Range("O87").Value2 = "17/08/2023 08:13"
$ws.Range("P87").Value2 = 4.99
$ws.Range("Q87").Value2 = "25/10/2023 13:56"
$ws.Range("R87").Value2 = 1.45
$ws.Range("S87").Value2 = "17/08/2023 08:13"
$ws.Range("T87").Value2 = 1.31
$ws.Range("U87").Value2 = "25/10/2023 13:44"
$ws.Range("V87").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-1923-partizan/l0enQa49/"
